# إضافة حدث جديد في Card2 by HOSSAM at 2025-12-08 11:50:59
#
# The author's source data (a pandas-style export) re-wrote the sheet with
# explicit "nan" placeholders for every still-empty data cell in the
# existing rows (2-13, columns D:O), and appended one new service-history
# row (row 14) recording a new maintenance event.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card2")

# --- Backfill empty data cells (rows 2-13, cols D:O) with literal "nan" ---
$ws.Range("D2:L2").Value = "nan"
$ws.Range("N2").Value = "nan"

$ws.Range("D3").Value = "nan"
$ws.Range("G3:K3").Value = "nan"
$ws.Range("M3:O3").Value = "nan"

$ws.Range("D4").Value = "nan"
$ws.Range("H4:K4").Value = "nan"
$ws.Range("M4:O4").Value = "nan"

$ws.Range("D5:O5").Value = "nan"

$ws.Range("E6:G6").Value = "nan"
$ws.Range("I6:K6").Value = "nan"
$ws.Range("M6:O6").Value = "nan"

$ws.Range("E7").Value = "nan"
$ws.Range("G7:J7").Value = "nan"
$ws.Range("M7:O7").Value = "nan"

$ws.Range("E8").Value = "nan"
$ws.Range("H8:K8").Value = "nan"
$ws.Range("M8:N8").Value = "nan"

$ws.Range("F9:K9").Value = "nan"
$ws.Range("N9").Value = "nan"

$ws.Range("D10:O10").Value = "nan"
$ws.Range("D11:O11").Value = "nan"
$ws.Range("D12:O12").Value = "nan"
$ws.Range("D13:O13").Value = "nan"

# --- Append the new service event as row 14 ---
# "card" column is stored as text ("2"), same as the existing rows above it,
# so force text formatting before assigning to avoid Excel's automatic
# number coercion.
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "2"
$ws.Range("L14").NumberFormat = "@"
$ws.Range("L14").Value = "30/6/2025"
$ws.Range("M14").Value = "قطع سير كويلر مسنن 1270"
$ws.Range("N14").Value = "تم تغير سير 1270"
$ws.Range("O14").Value = "فني"
